$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A width (closest achievable approximation of 15.42578125 chars width)
$ws.Columns.Item(1).ColumnWidth = 14.59

# Update values in column A (rows 1-33)
$values = @(
    0.072771377589418762,
    -0.0059999999515376601,
    -0.003999999962811529,
    -0.0079999999305968572,
    -0.0029999999702221558,
    -0.0019999999787660983,
    -0.0099999999149016361,
    -0.0099999999135089723,
    -0.0019999999766007193,
    -0.0019999999765385468,
    -0.0029999999686172174,
    -0.0034999999647213897,
    -0.0034999999664488968,
    -0.0079999999317195147,
    0.01426016847050704,
    -0.0019999999814341862,
    -0.0019999999826838533,
    0.0078654463508174288,
    -0.0039999999672644115,
    0.0077953767099359084,
    -0.0039999999664974695,
    -0.0039999999660853547,
    -0.0049999999521457283,
    -0.019999999827506443,
    -0.019999999824698911,
    -0.0024999999755337399,
    0.072416174211952899,
    -0.0019999999737327911,
    -0.0069999999289400705,
    -0.059999999506005874,
    -0.0069999999313772321,
    0.039924424325612051,
    -0.0039999999568038902
)

For ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
